$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.564.21"
$ws.Range("E2").Value = "  +6.88%  "
$ws.Range("D3").Value = "3.326.15"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "410.40"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").Value = "115.76"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "3.321.66"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D10").Value = "0.627"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.114"
$ws.Range("E11").Value = "  +17.50%  "
$ws.Range("D12").Value = "40.11"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "3.845.00"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "8.22"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "19.19"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "3.327.78"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "60.468.43"
$ws.Range("E18").Value = "  +6.66%  "
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "10.83"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "3.37"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0000114"
$ws.Range("E22").Value = "  +5.06%  "
$ws.Range("D23").Value = "12.41"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").Value = "295.79"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "73.88"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").Value = "3.12"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "29.11"
$ws.Range("E27").Value = "  +3.61%  "
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").Value = "  +5.14%  "
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "7.51"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "0.113"
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "11.31"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "41.27"
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").Value = "2.48"
$ws.Range("E36").Value = "  +16.42%  "
$ws.Range("D37").Value = "0.0494"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "52.05"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "3.05"
$ws.Range("E40").Value = "  +5.82%  "
$ws.Range("D41").Value = "3.37"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").Value = "133.79"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "3.88"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "16.28"
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "21.17"
$ws.Range("E49").Value = "  -4.66%  "
$ws.Range("D50").Value = "2.139.45"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "3.649.16"
$ws.Range("E51").Value = "  +2.10%  "
